$wb = $excel.ActiveWorkbook

$sourceTypeList = '"repository_id,file,restriction,ligation,PCR,homologous_recombination,gibson_assembly,restriction_and_ligation,genome_coordinates,manually_typed"'

# --- Sheet 1: NamedThing -- drop "name" and "description" columns, keep only "id" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1:C1").ClearContents()

# --- Sheet 2: CloningStrategy -> ManuallyTypedSource ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ManuallyTypedSource"
$ws2.Range("D2:D1048576").Validation.Delete()
$ws2.Cells.ClearContents()
$ws2.Range("A1").Value = "user_input"
$ws2.Range("B1").Value = "input"
$ws2.Range("C1").Value = "output"
$ws2.Range("D1").Value = "type"
$ws2.Range("E1").Value = "kind"
$ws2.Range("F1").Value = "info"
$ws2.Range("G1").Value = "id"
$ws2.Range("D2:D1048576").Validation.Add(3, 1, 1, $sourceTypeList)

# --- Sheet 3: CloningStrategyCollection -> UploadedFileSource ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "UploadedFileSource"
$ws3.Cells.ClearContents()
$ws3.Range("A1").Value = "file_name"
$ws3.Range("B1").Value = "index_in_file"
$ws3.Range("C1").Value = "input"
$ws3.Range("D1").Value = "output"
$ws3.Range("E1").Value = "type"
$ws3.Range("F1").Value = "kind"
$ws3.Range("G1").Value = "info"
$ws3.Range("H1").Value = "id"
$ws3.Range("E2:E1048576").Validation.Add(3, 1, 1, $sourceTypeList)

# --- New Sheet 4: RepositoryIdSource ---
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "RepositoryIdSource"
$ws4.Range("A1").Value = "repository_name"
$ws4.Range("B1").Value = "repository_id"
$ws4.Range("C1").Value = "input"
$ws4.Range("D1").Value = "output"
$ws4.Range("E1").Value = "type"
$ws4.Range("F1").Value = "kind"
$ws4.Range("G1").Value = "info"
$ws4.Range("H1").Value = "id"
$ws4.Range("A2:A1048576").Validation.Add(3, 1, 1, '"addgene,genbank"')
$ws4.Range("E2:E1048576").Validation.Add(3, 1, 1, $sourceTypeList)

# --- New Sheet 5: GenomeCoordinatesSource ---
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "GenomeCoordinatesSource"
$ws5.Range("A1").Value = "assembly_accession"
$ws5.Range("B1").Value = "sequence_accession"
$ws5.Range("C1").Value = "locus_tag"
$ws5.Range("D1").Value = "gene_id"
$ws5.Range("E1").Value = "start"
$ws5.Range("F1").Value = "stop"
$ws5.Range("G1").Value = "strand"
$ws5.Range("H1").Value = "input"
$ws5.Range("I1").Value = "output"
$ws5.Range("J1").Value = "type"
$ws5.Range("K1").Value = "kind"
$ws5.Range("L1").Value = "info"
$ws5.Range("M1").Value = "id"
$ws5.Range("J2:J1048576").Validation.Add(3, 1, 1, $sourceTypeList)
